$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("D2").Value = -0.251

$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = -2.03
$ws.Range("L2").Value = -4.194214876033057
$ws.Range("M2").Value = 0.012
$ws.Range("N2").Value = 0.003468208092485549
$ws.Range("O2").Value = -0.005911330049261084
$ws.Range("P2").Value = 0.012
$ws.Range("Q2").Value = 0.003468208092485549
$ws.Range("R2").Value = -0.005911330049261084

$ws.Range("U2").Value = 4.32
$ws.Range("V2").Value = 1.248554913294798
$ws.Range("W2").Value = -0.179646017699115
$ws.Range("X2").Value = 1.129256214249508
$ws.Range("Y2").Value = -1.308902231948623
$ws.Range("Z2").Value = 0.004337694927406345

$ws.Range("AB2").Value = 0.09795707071344843
$ws.Range("AC2").Value = -0.09795707071344843
$ws.Range("AD2").Value = 91.5
$ws.Range("AF2").Value = 91.5
$ws.Range("AG2").Value = 87.18000000000001
$ws.Range("AH2").Value = 0.9635636057287279
$ws.Range("AI2").Value = 0.9072880515617254
$ws.Range("AJ2").Value = 0.9618270079435128
$ws.Range("AK2").Value = 0.9031389205428365

# --- Row 3 updates ---
$ws.Range("B3").Value = "BH Leasing Société Anonyme (BVMT:BHL)"
$ws.Range("D3").Value = -0.251

$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = -2.03
$ws.Range("L3").Value = -4.194214876033057
$ws.Range("M3").Value = 0.012
$ws.Range("N3").Value = 0.003468208092485549
$ws.Range("O3").Value = -0.005911330049261084
$ws.Range("P3").Value = 0.012
$ws.Range("Q3").Value = 0.003468208092485549
$ws.Range("R3").Value = -0.005911330049261084

$ws.Range("U3").Value = 4.32
$ws.Range("V3").Value = 1.248554913294798
$ws.Range("W3").Value = -0.179646017699115
$ws.Range("X3").Value = 1.129256214249508
$ws.Range("Y3").Value = -1.308902231948623
$ws.Range("Z3").Value = 0.004337694927406345

$ws.Range("AB3").Value = 0.09795707071344843
$ws.Range("AC3").Value = -0.09795707071344843
$ws.Range("AD3").Value = 91.5
$ws.Range("AF3").Value = 91.5
$ws.Range("AG3").Value = 87.18000000000001
$ws.Range("AH3").Value = 0.9635636057287279
$ws.Range("AI3").Value = 0.9072880515617254
$ws.Range("AJ3").Value = 0.9618270079435128
$ws.Range("AK3").Value = 0.9031389205428365

$wb.Save()
